$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.218.53"
$ws.Range("E2").Value = "'  +0.70%  "
$ws.Range("D3").Value = "'1.853.92"
$ws.Range("E3").Value = "'  +1.31%  "
$ws.Range("E4").Value = "'  -0.40%  "
$ws.Range("D5").Value = "'313.93"
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "'  -0.31%  "
$ws.Range("D7").Value = "'0.4634"
$ws.Range("E7").Value = "'  +0.12%  "
$ws.Range("D8").Value = "'0.3712"
$ws.Range("E8").Value = "'  -0.01%  "
$ws.Range("D9").Value = "'0.07287"
$ws.Range("E9").Value = "'  -0.79%  "
$ws.Range("D10").Value = "'0.8866"
$ws.Range("E10").Value = "'  +1.00%  "
$ws.Range("D11").Value = "'20.11"
$ws.Range("E11").Value = "'  +1.44%  "
$ws.Range("D12").Value = "'0.07869"
$ws.Range("E12").Value = "'  -0.20%  "
$ws.Range("D13").Value = "'1.899.20"
$ws.Range("E13").Value = "'  +4.60%  "
$ws.Range("E14").Value = "'  +0.91%  "
$ws.Range("D15").Value = "'6.515"
$ws.Range("E15").Value = "'  -0.57%  "
$ws.Range("D16").Value = "'91.05"
$ws.Range("E16").Value = "'  -0.53%  "
$ws.Range("D18").Value = "'0.000008922"
$ws.Range("E18").Value = "'  +0.84%  "
$ws.Range("E19").Value = "'  -0.25%  "
$ws.Range("D20").Value = "'14.68"
$ws.Range("E20").Value = "'  -0.79%  "
$ws.Range("D21").Value = "'27.250.25"
$ws.Range("E21").Value = "'  +0.76%  "
$ws.Range("D22").Value = "'5.084"
$ws.Range("E22").Value = "'  -0.47%  "
$ws.Range("E23").Value = "'  -0.18%  "
$ws.Range("D24").Value = "'2.097.30"
$ws.Range("E24").Value = "'  +3.44%  "
$ws.Range("D25").Value = "'1.950"
$ws.Range("E25").Value = "'  +5.23%  "
$ws.Range("D26").Value = "'151.30"
$ws.Range("E26").Value = "'  -0.98%  "
$ws.Range("D27").Value = "'18.39"
$ws.Range("E27").Value = "'  -0.55%  "
$ws.Range("D28").Value = "'2.049"
$ws.Range("E28").Value = "'  -0.01%  "
$ws.Range("D29").Value = "'115.96"
$ws.Range("D30").Value = "'5.043"
$ws.Range("E30").Value = "'  -1.73%  "
$ws.Range("D31").Value = "'0.08803"
$ws.Range("E31").Value = "'  -1.00%  "
$ws.Range("D32").Value = "'3.140"
$ws.Range("E32").Value = "'  +6.19%  "
$ws.Range("D33").Value = "'0.7697"
$ws.Range("E33").Value = "'  +5.51%  "
$ws.Range("D34").Value = "'1.165"
$ws.Range("E34").Value = "'  +2.76%  "
$ws.Range("D35").Value = "'4.520"
$ws.Range("E35").Value = "'  +1.84%  "
$ws.Range("D36").Value = "'2.729"
$ws.Range("E36").Value = "'  +10.79%  "
$ws.Range("D37").Value = "'1.106"
$ws.Range("E37").Value = "'  +2.57%  "
$ws.Range("D38").Value = "'0.01939"
$ws.Range("E38").Value = "'  -0.68%  "
$ws.Range("D39").Value = "'0.05219"
$ws.Range("E39").Value = "'  -0.26%  "
$ws.Range("D40").Value = "'2.941"
$ws.Range("E40").Value = "'  -0.35%  "
$ws.Range("D41").Value = "'7.034"
$ws.Range("E41").Value = "'  -1.16%  "
$ws.Range("D42").Value = "'0.5125"
$ws.Range("E42").Value = "'  -0.93%  "
$ws.Range("D43").Value = "'0.1629"
$ws.Range("E43").Value = "'  +0.16%  "
$ws.Range("D44").Value = "'8.456"
$ws.Range("E44").Value = "'  +3.29%  "
$ws.Range("D45").Value = "'0.4799"
$ws.Range("E45").Value = "'  -0.84%  "
$ws.Range("D46").Value = "'10.38"
$ws.Range("E46").Value = "'  +1.98%  "
$ws.Range("D47").Value = "'1.001"
$ws.Range("E47").Value = "'  -0.35%  "
$ws.Range("D48").Value = "'102.63"
$ws.Range("E48").Value = "'  +0.21%  "
$ws.Range("D49").Value = "'1.644"
$ws.Range("E49").Value = "'  +0.68%  "
$ws.Range("D50").Value = "'0.06202"
$ws.Range("E50").Value = "'  -0.04%  "
$ws.Range("D51").Value = "'65.31"
$ws.Range("E51").Value = "'  +0.85%  "
